# Fixed #450 Add support for text-decoration:line-through style in
# fromHTMLBodyString() service.
#
# The HTML->Word converter now always stamps the run's bold / italic /
# strike-through state explicitly (b / i / strike) instead of leaving
# the run properties implicit, so that "no bold / no italic / no
# strike-through" coming out of the HTML source round-trips instead of
# silently inheriting whatever the destination style defines. Apply
# that to the REC/RPL bullet list runs (and the hyperlink runs amongst
# them) in the "Network" section.

$d = $word.ActiveDocument
$nbsp = [char]0x00A0

function Stamp-NoFormatting($range) {
    $range.Font.Bold = $false
    $range.Font.Italic = $false
    $range.Font.StrikeThrough = $false
}

function Find-Exact($needle) {
    $r = $d.Content
    $ok = $r.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "text not found: $needle"
    }
    return $r
}

# 1) "REC - Unit Network Adapter" hyperlink run
$r = Find-Exact("REC - Unit Network Adapter")
Stamp-NoFormatting $r

# 2) NBSP + "describes the REC (i.e.what will be replicated) " run
$r = Find-Exact($nbsp + "describes the REC (i.e.what will be replicated) ")
Stamp-NoFormatting $r

# 3) "RPL - Instantiations of Unit Network Adapter" hyperlink run
$r = Find-Exact("RPL - Instantiations of Unit Network Adapter")
Stamp-NoFormatting $r

# 4) NBSP + "describes the three RPLs (replicas). " run
$r = Find-Exact($nbsp + "describes the three RPLs (replicas). ")
Stamp-NoFormatting $r

# 5) " The connection between RPLs" + NBSP + "is detailed in " run
$r = Find-Exact(" The connection between RPLs" + $nbsp + "is detailed in ")
Stamp-NoFormatting $r

# 6) "[PAB] Focus on Network Setup, Configuration and Tests" hyperlink run
$r = Find-Exact("[PAB] Focus on Network Setup, Configuration and Tests")
Stamp-NoFormatting $r
